$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("data")

# Update the two timestamp cells on the "data" sheet
$ws1.Range("F2").Value = "2021-10-05 14:19:35.749733"
$ws1.Range("F3").Value = "2021-10-05 14:19:35.749741"

# Add a new "metadata" worksheet right after "data"
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "metadata"

# Reuse the same header style (bold + border + center/top alignment) used
# on the "data" sheet's header row, and the same style used for the
# numeric index column, so no new style table entries are introduced.
$ws1.Range("B1:F1").Copy()
$ws2.Range("B1:F1").PasteSpecial(-4122)
$ws1.Range("B1").Copy()
$ws2.Range("G1").PasteSpecial(-4122)

$ws1.Range("A2").Copy()
$ws2.Range("A2").PasteSpecial(-4122)

# Header row
$ws2.Range("B1").Value = "data_name"
$ws2.Range("C1").Value = "data_id"
$ws2.Range("D1").Value = "data_version"
$ws2.Range("E1").Value = "data_version_created"
$ws2.Range("F1").Value = "panel_query_time"
$ws2.Range("G1").Value = "panel_get_request"

# Data row
$ws2.Range("A2").Value = 0
$ws2.Range("B2").Value = "Classical tuberous sclerosis"
$ws2.Range("C2").Value = 197

# "1.2" must remain text, not be coerced into the number 1.2
$verCell = $ws2.Cells.Item(2, 4)
$verCell.NumberFormat = "@"
$verCell.Value = "1.2"
$verCell.Style = "Normal"

$ws2.Range("E2").Value = "2017-11-05T02:37:20.208587Z"
$ws2.Range("F2").Value = "2021-10-05 14:19:35.746595"
$ws2.Range("G2").Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/197/?format=json"

Write-Output "metadata sheet added"
